$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$f = $sec.Footers.Item(1)
$r = $f.Range
Write-Output ("Footer paragraphs count: " + $r.Paragraphs.Count)
$p2 = $r.Paragraphs.Item(3)
Write-Output ("p2(3) text=[" + $p2.Range.Text + "] start=" + $p2.Range.Start + " end=" + $p2.Range.End)
